$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1.459612070389937, 1.667794583268128, 0.1575252929769615, 0.496779210170732),
    @(0.3048080303191223, 1.667794583268128, 26.21740644021617, 8.660232485948974),
    @(3.230985683306322, 1.667794583268128, 0.8054896365839992, 0.496779210170732),
    @(0.127881588408715, 0.3127903958511391, 3.900430680208489, 8.660232485948974),
    @(0.3048080303191223, 0.04240448674262143, 0.8054896365839992, 0.496779210170732),
    @(0.6753301551942219, 1.667794583268128, 0.1575252929769615, 0.496779210170732),
    @(3.230985683306322, 1.667794583268128, 0.1575252929769615, 8.660232485948974),
    @(3.230985683306322, 0.3127903958511391, 0.1575252929769615, 0.496779210170732)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $b = $vals[0]
    $c = $vals[1]
    $d = $vals[2]
    $e = $vals[3]
    $g = $b + $c + $d + $e

    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 7).Value = $g
}
